# Auto-generated: update leve profit market-data values across sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 20624.777  # H15
$ws.Cells.Item(15, 9).Value = 20624.777  # I15
$ws.Cells.Item(15, 11).Value = 61874.33099999999  # K15
$ws.Cells.Item(15, 13).Value = -61705.33099999999  # M15
$ws.Cells.Item(61, 8).Value = 99  # H61
$ws.Cells.Item(61, 9).Value = 99  # I61
$ws.Cells.Item(61, 11).Value = 297  # K61
$ws.Cells.Item(61, 13).Value = -125  # M61
$ws.Cells.Item(101, 8).Value = 899.5  # H101
$ws.Cells.Item(101, 10).Value = 1328  # J101
$ws.Cells.Item(101, 12).Value = 3984  # L101
$ws.Cells.Item(101, 14).Value = -7228  # N101
$ws.Cells.Item(104, 8).Value = 1147.25  # H104
$ws.Cells.Item(104, 9).Value = 1147.25  # I104
$ws.Cells.Item(104, 11).Value = 3441.75  # K104
$ws.Cells.Item(104, 13).Value = -1694.75  # M104
$ws.Cells.Item(112, 8).Value = 4839.185  # H112
$ws.Cells.Item(112, 10).Value = 4839.185  # J112
$ws.Cells.Item(112, 12).Value = 14517.555  # L112
$ws.Cells.Item(112, 14).Value = -16733.555  # N112
$ws.Cells.Item(115, 8).Value = 0  # H115
$ws.Cells.Item(115, 9).Value = 0  # I115
$ws.Cells.Item(115, 11).Value = 0  # K115
$ws.Cells.Item(115, 13).Value = $null  # M115
$ws.Cells.Item(118, 8).Value = 511.14285  # H118
$ws.Cells.Item(118, 9).Value = 521.3333  # I118
$ws.Cells.Item(118, 11).Value = 1563.9999  # K118
$ws.Cells.Item(118, 13).Value = 93.00009999999997  # M118
$ws.Cells.Item(127, 8).Value = 773.2  # H127
$ws.Cells.Item(127, 9).Value = 941.5  # I127
$ws.Cells.Item(127, 10).Value = 100  # J127
$ws.Cells.Item(127, 11).Value = 2824.5  # K127
$ws.Cells.Item(127, 12).Value = 300  # L127
$ws.Cells.Item(127, 13).Value = 2135.5  # M127
$ws.Cells.Item(127, 14).Value = -10220  # N127
$ws.Cells.Item(129, 8).Value = 1329.2727  # H129
$ws.Cells.Item(129, 9).Value = 770.3333  # I129
$ws.Cells.Item(129, 11).Value = 2310.9999  # K129
$ws.Cells.Item(129, 13).Value = 2689.0001  # M129
$ws.Cells.Item(137, 8).Value = 5883.354  # H137
$ws.Cells.Item(137, 9).Value = 4212.1816  # I137
$ws.Cells.Item(137, 10).Value = 9559.933999999999  # J137
$ws.Cells.Item(137, 11).Value = 12636.5448  # K137
$ws.Cells.Item(137, 12).Value = 28679.802  # L137
$ws.Cells.Item(137, 13).Value = -10086.5448  # M137
$ws.Cells.Item(137, 14).Value = -33779.802  # N137
$ws.Cells.Item(138, 8).Value = 3148.29  # H138
$ws.Cells.Item(138, 10).Value = 3233.9285  # J138
$ws.Cells.Item(138, 12).Value = 9701.7855  # L138
$ws.Cells.Item(138, 14).Value = -19981.7855  # N138

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(19, 8).Value = 499.5  # H19
$ws.Cells.Item(19, 9).Value = 490  # I19
$ws.Cells.Item(19, 10).Value = 509  # J19
$ws.Cells.Item(19, 11).Value = 490  # K19
$ws.Cells.Item(19, 12).Value = 509  # L19
$ws.Cells.Item(19, 13).Value = -261  # M19
$ws.Cells.Item(19, 14).Value = -967  # N19
$ws.Cells.Item(32, 8).Value = 2432.31  # H32
$ws.Cells.Item(32, 9).Value = 2219.969  # I32
$ws.Cells.Item(32, 11).Value = 2219.969  # K32
$ws.Cells.Item(32, 13).Value = -1932.969  # M32
$ws.Cells.Item(61, 8).Value = 22227728  # H61
$ws.Cells.Item(61, 9).Value = 2447.2  # I61
$ws.Cells.Item(61, 11).Value = 2447.2  # K61
$ws.Cells.Item(61, 13).Value = -2235.2  # M61
$ws.Cells.Item(74, 8).Value = 20470.38  # H74
$ws.Cells.Item(74, 9).Value = 26198.176  # I74
$ws.Cells.Item(74, 10).Value = 5196.2666  # J74
$ws.Cells.Item(74, 11).Value = 26198.176  # K74
$ws.Cells.Item(74, 12).Value = 5196.2666  # L74
$ws.Cells.Item(74, 13).Value = -25324.176  # M74
$ws.Cells.Item(74, 14).Value = -6944.2666  # N74
$ws.Cells.Item(77, 8).Value = 20470.38  # H77
$ws.Cells.Item(77, 9).Value = 26198.176  # I77
$ws.Cells.Item(77, 10).Value = 5196.2666  # J77
$ws.Cells.Item(77, 11).Value = 130990.88  # K77
$ws.Cells.Item(77, 12).Value = 25981.333  # L77
$ws.Cells.Item(77, 13).Value = -126622.88  # M77
$ws.Cells.Item(77, 14).Value = -34717.333  # N77
$ws.Cells.Item(102, 8).Value = 5001963  # H102
$ws.Cells.Item(102, 9).Value = 5407242.5  # I102
$ws.Cells.Item(102, 11).Value = 5407242.5  # K102
$ws.Cells.Item(102, 13).Value = -5405620.5  # M102
$ws.Cells.Item(136, 8).Value = 22227728  # H136
$ws.Cells.Item(136, 9).Value = 2447.2  # I136
$ws.Cells.Item(136, 11).Value = 7341.599999999999  # K136
$ws.Cells.Item(136, 13).Value = -4791.599999999999  # M136

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(106, 8).Value = 23142.428  # H106
$ws.Cells.Item(106, 10).Value = 23142.428  # J106
$ws.Cells.Item(106, 12).Value = 23142.428  # L106
$ws.Cells.Item(106, 14).Value = -25666.428  # N106
$ws.Cells.Item(113, 8).Value = 5140.7144  # H113
$ws.Cells.Item(113, 9).Value = 5140.7144  # I113
$ws.Cells.Item(113, 11).Value = 5140.7144  # K113
$ws.Cells.Item(113, 13).Value = -2970.7144  # M113
$ws.Cells.Item(134, 8).Value = 4550352.5  # H134
$ws.Cells.Item(134, 9).Value = 6100112  # I134
$ws.Cells.Item(134, 11).Value = 18300336  # K134
$ws.Cells.Item(134, 13).Value = -18297801  # M134

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2587  # H16
$ws.Cells.Item(16, 9).Value = 1624.5  # I16
$ws.Cells.Item(16, 11).Value = 1624.5  # K16
$ws.Cells.Item(16, 13).Value = -1337.5  # M16
$ws.Cells.Item(43, 8).Value = 30045  # H43
$ws.Cells.Item(43, 10).Value = 30045  # J43
$ws.Cells.Item(43, 12).Value = 30045  # L43
$ws.Cells.Item(43, 14).Value = -30413  # N43
$ws.Cells.Item(101, 8).Value = 30045  # H101
$ws.Cells.Item(101, 10).Value = 30045  # J101
$ws.Cells.Item(101, 12).Value = 30045  # L101
$ws.Cells.Item(101, 14).Value = -36535  # N101
$ws.Cells.Item(113, 8).Value = 2587  # H113
$ws.Cells.Item(113, 9).Value = 1624.5  # I113
$ws.Cells.Item(113, 11).Value = 1624.5  # K113
$ws.Cells.Item(113, 13).Value = 545.5  # M113
$ws.Cells.Item(141, 8).Value = 238865.47  # H141
$ws.Cells.Item(141, 10).Value = 263380.88  # J141
$ws.Cells.Item(141, 12).Value = 263380.88  # L141
$ws.Cells.Item(141, 14).Value = -273740.88  # N141

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(76, 8).Value = 3901.5  # H76
$ws.Cells.Item(76, 9).Value = 3901.5  # I76
$ws.Cells.Item(76, 11).Value = 11704.5  # K76
$ws.Cells.Item(76, 13).Value = -11321.5  # M76
$ws.Cells.Item(79, 8).Value = 3901.5  # H79
$ws.Cells.Item(79, 9).Value = 3901.5  # I79
$ws.Cells.Item(79, 11).Value = 11704.5  # K79
$ws.Cells.Item(79, 13).Value = -10378.5  # M79
$ws.Cells.Item(113, 8).Value = 1662.8  # H113
$ws.Cells.Item(113, 10).Value = 1852.7646  # J113
$ws.Cells.Item(113, 12).Value = 5558.293799999999  # L113
$ws.Cells.Item(113, 14).Value = -9898.293799999999  # N113
$ws.Cells.Item(122, 8).Value = 3192506.2  # H122
$ws.Cells.Item(122, 9).Value = 7073333  # I122
$ws.Cells.Item(122, 10).Value = 1252093  # J122
$ws.Cells.Item(122, 11).Value = 63659997  # K122
$ws.Cells.Item(122, 12).Value = 11268837  # L122
$ws.Cells.Item(122, 13).Value = -63657547  # M122
$ws.Cells.Item(122, 14).Value = -11273737  # N122

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(52, 8).Value = 0  # H52
$ws.Cells.Item(52, 10).Value = 0  # J52
$ws.Cells.Item(52, 12).Value = 0  # L52
$ws.Cells.Item(52, 14).Value = $null  # N52
$ws.Cells.Item(93, 8).Value = 39985  # H93
$ws.Cells.Item(93, 9).Value = 0  # I93
$ws.Cells.Item(93, 10).Value = 39985  # J93
$ws.Cells.Item(93, 11).Value = 0  # K93
$ws.Cells.Item(93, 12).Value = 39985  # L93
$ws.Cells.Item(93, 13).Value = $null  # M93
$ws.Cells.Item(93, 14).Value = -43729  # N93

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(13, 8).Value = 3996.6667  # H13
$ws.Cells.Item(13, 9).Value = 4990  # I13
$ws.Cells.Item(13, 11).Value = 4990  # K13
$ws.Cells.Item(13, 13).Value = -4850  # M13
$ws.Cells.Item(17, 8).Value = 1777.4166  # H17
$ws.Cells.Item(17, 9).Value = 1777.4166  # I17
$ws.Cells.Item(17, 11).Value = 1777.4166  # K17
$ws.Cells.Item(17, 13).Value = -1607.4166  # M17
$ws.Cells.Item(100, 8).Value = 7146.154  # H100
$ws.Cells.Item(100, 9).Value = 6387  # I100
$ws.Cells.Item(100, 11).Value = 6387  # K100
$ws.Cells.Item(100, 13).Value = -5846  # M100
$ws.Cells.Item(103, 8).Value = 25000  # H103
$ws.Cells.Item(103, 10).Value = 25000  # J103
$ws.Cells.Item(103, 12).Value = 25000  # L103
$ws.Cells.Item(103, 14).Value = -27344  # N103
$ws.Cells.Item(122, 8).Value = 4026.7222  # H122
$ws.Cells.Item(122, 9).Value = 3444.95  # I122
$ws.Cells.Item(122, 10).Value = 5688.9287  # J122
$ws.Cells.Item(122, 11).Value = 10334.85  # K122
$ws.Cells.Item(122, 12).Value = 17066.7861  # L122
$ws.Cells.Item(122, 13).Value = -7884.849999999999  # M122
$ws.Cells.Item(122, 14).Value = -21966.7861  # N122
$ws.Cells.Item(136, 8).Value = 13057.521  # H136
$ws.Cells.Item(136, 9).Value = 3742.5715  # I136
$ws.Cells.Item(136, 10).Value = 14729.436  # J136
$ws.Cells.Item(136, 11).Value = 11227.7145  # K136
$ws.Cells.Item(136, 12).Value = 44188.308  # L136
$ws.Cells.Item(136, 13).Value = -8677.7145  # M136
$ws.Cells.Item(136, 14).Value = -49288.308  # N136

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(69, 8).Value = 43429.5  # H69
$ws.Cells.Item(69, 10).Value = 44063  # J69
$ws.Cells.Item(69, 12).Value = 44063  # L69
$ws.Cells.Item(69, 14).Value = -45561  # N69
$ws.Cells.Item(72, 8).Value = 43429.5  # H72
$ws.Cells.Item(72, 10).Value = 44063  # J72
$ws.Cells.Item(72, 12).Value = 132189  # L72
$ws.Cells.Item(72, 14).Value = -139677  # N72
$ws.Cells.Item(136, 8).Value = 10102937  # H136
$ws.Cells.Item(136, 9).Value = 13159630  # I136
$ws.Cells.Item(136, 11).Value = 39478890  # K136
$ws.Cells.Item(136, 13).Value = -39476340  # M136

